$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 112, shifting existing rows 112-178 down to 113-179
$ws.Rows.Item(112).Insert()

# Fill in the new row 112 with data (same template values as surrounding rows,
# with the new record's specific values)
$ws.Cells.Item(112, 1).Value = 4
$ws.Cells.Item(112, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(112, 3).Value = "Los Lagos"
$ws.Cells.Item(112, 4).Value = 45001
$ws.Cells.Item(112, 5).Value = 10
$ws.Cells.Item(112, 6).Value = 100112052
$ws.Cells.Item(112, 7).Value = "Albahaca"
$ws.Cells.Item(112, 8).Value = "Sin especificar"
$ws.Cells.Item(112, 9).Value = "Primera"
$ws.Cells.Item(112, 10).Value = 60
$ws.Cells.Item(112, 11).Value = 7000
$ws.Cells.Item(112, 12).Value = 7000
$ws.Cells.Item(112, 13).Value = 7000
$ws.Cells.Item(112, 14).Value = "`$/docena de matas"
$ws.Cells.Item(112, 15).Value = "Región Metropolitana"
$ws.Cells.Item(112, 16).Value = 1167
$ws.Cells.Item(112, 17).Value = 6
$ws.Cells.Item(112, 18).Value = "Hortaliza"
